$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the existing "sum" header cell (G1) onto the new
# "Save" header cell (H1) so it reuses the same cell style as the other
# header cells.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# Header cell for the new "Save" column
$ws.Range("H1").Value = "Save"

# Data values for the new column
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 1
$ws.Range("H4").Value = 1
